# Updates crypto price/volume data to the latest snapshot (GitHub Actions scrape).
# Column D ("Price") must stay plain text (matches source workbook, which stores
# every price as inlineStr, e.g. "67.509.98"), so we force a text NumberFormat
# before writing then restore the default "Normal" style so no spurious style
# index is left on the cell. Column E ("Volume(1h)") already round-trips as text
# because of the surrounding spaces/"%" sign, so it needs no special handling.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "67.509.98"
$ws.Range("D2").Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +0.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.527.95"
$ws.Range("D3").Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -0.17%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "615.24"
$ws.Range("D5").Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "152.01"
$ws.Range("D6").Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -1.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "3.527.63"
$ws.Range("D7").Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -0.10%  "

$ws.Cells.Item(8, 5).Value = "  -0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.481"
$ws.Range("D9").Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -0.88%  "

$ws.Cells.Item(10, 5).Value = "  -0.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "7.10"
$ws.Range("D11").Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +2.97%  "

$ws.Cells.Item(12, 5).Value = "  -1.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.0000222"
$ws.Range("D13").Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -0.63%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "32.19"
$ws.Range("D14").Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -0.23%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "4.123.39"
$ws.Range("D15").Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +0.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "3.525.68"
$ws.Range("D16").Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -0.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "67.466.31"
$ws.Range("D17").Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -0.13%  "

$ws.Cells.Item(18, 5).Value = "  -0.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "6.41"
$ws.Range("D19").Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +0.80%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "15.35"
$ws.Range("D20").Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -1.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "445.68"
$ws.Range("D21").Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -1.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "9.52"
$ws.Range("D22").Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +1.09%  "

$ws.Cells.Item(23, 5).Value = "  -2.89%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "77.47"
$ws.Range("D24").Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -0.94%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.0000132"
$ws.Range("D25").Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +10.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "3.667.49"
$ws.Range("D26").Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -0.30%  "

$ws.Cells.Item(27, 5).Value = "  +0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "10.27"
$ws.Range("D28").Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -2.45%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "8.53"
$ws.Range("D29").Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +1.94%  "

$ws.Cells.Item(30, 5).Value = "  -2.49%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.56"
$ws.Range("D31").Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -7.46%  "

$ws.Cells.Item(32, 5).Value = "  -0.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.164"
$ws.Range("D33").Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +4.21%  "

$ws.Cells.Item(34, 5).Value = "  -0.35%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "6.19"
$ws.Range("D35").Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -0.66%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "3.519.50"
$ws.Range("D36").Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -0.46%  "

$ws.Cells.Item(37, 5).Value = "  -3.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "8.04"
$ws.Range("D38").Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +0.36%  "

$ws.Cells.Item(39, 5).Value = "  +0.00%  "

$ws.Cells.Item(40, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -0.08%  "

$ws.Cells.Item(41, 2).Value = "Monero"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "177.28"
$ws.Range("D41").Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +0.51%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "2.17"
$ws.Range("D42").Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +3.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.0884"
$ws.Range("D43").Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +0.60%  "

$ws.Cells.Item(44, 5).Value = "  -2.94%  "

$ws.Cells.Item(45, 5).Value = "  -0.65%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "28.54"
$ws.Range("D46").Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -2.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "45.16"
$ws.Range("D47").Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -1.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "2.64"
$ws.Range("D48").Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +2.13%  "

$ws.Cells.Item(49, 5).Value = "  +3.67%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "7.61"
$ws.Range("D50").Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -0.68%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.00"
$ws.Range("D51").Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -1.22%  "
